$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy number formatting from the (now-shifted) column E into the new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest period values
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 5207000
$ws.Range("D9").Value = 3181000
$ws.Range("D10").Value = 2026000
$ws.Range("D12").Value = 189000
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 22000
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 4553000
$ws.Range("D18").Value = 654000
$ws.Range("D20").Value = 13000
$ws.Range("D21").Value = 928000
$ws.Range("D22").Value = 82000
$ws.Range("D23").Value = 585000
$ws.Range("D24").Value = 28500
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 556500
$ws.Range("D27").Value = 556500
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -7500
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -13000
$ws.Range("D33").Value = 549000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 549000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 296000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 1031000
$ws.Range("D44").Value = 595000
$ws.Range("D45").Value = 172000
$ws.Range("D46").Value = 2094000
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 656000
$ws.Range("D49").Value = 4208000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 264000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 7222000
$ws.Range("D57").Value = 586000
$ws.Range("D58").Value = 257000
$ws.Range("D59").Value = 546000
$ws.Range("D60").Value = 1389000
$ws.Range("D61").Value = 2051000
$ws.Range("D62").Value = 1000000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 4454000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 1639000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2768000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 549000
$ws.Range("D83").Value = 261000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 586000
$ws.Range("D91").Value = -237000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -643000
$ws.Range("D96").Value = -152000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -40000
$ws.Range("D101").Value = -21000
$ws.Range("D102").Value = -118000
